$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04) for every data row (rows 2-158).
$ws.Range("C2:C158").Value2 = 45203
